$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string used by F2/F3 ("pcl01@bpssumsel.com" -> "idris@bps.go.id")
$ws.Range("F2").Value = "idris@bps.go.id"
$ws.Range("F3").Value = "idris@bps.go.id"

# Remove the extra fill style (cellXfs index 4, applyFill) that was applied to F2/F3
# by clearing their formatting back to the default style (same as e.g. A2/C2 cells)
$ws.Range("F2:F3").Style = "Normal"

# Set new column F width (~21.5703125 "characters" as stored in the xlsx;
# the COM ColumnWidth property here rounds to the nearest 1/6 of a
# character, so 20.7 is the closest settable value that serializes back
# out to the target stored width of 21.5)
$ws.Columns.Item(6).ColumnWidth = 20.7

# Update the active selection cell to F7
$ws.Range("F7").Select()
